$wb = $excel.ActiveWorkbook

# Add the new "tmp" worksheet at the end (after the current last sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "tmp"

[void]$newSheet.Select()

# --- Header row for the combined sequence/group/GRAVY/charge table ---
$newSheet.Cells.Item(1, 1).Value = "sequence"
$newSheet.Cells.Item(1, 2).Value = "group"
$newSheet.Cells.Item(1, 3).Value = "GRAVY index"
$newSheet.Cells.Item(1, 4).Value = "Net charge at pH 7"

# --- "Ob" group source block (columns G:J), written before the Ob rows of the
#     combined table so shared strings intern in the same order as the source data ---
$newSheet.Cells.Item(2, 7).Value = "sequence"
$newSheet.Cells.Item(2, 8).Value = "GRAVY index"
$newSheet.Cells.Item(2, 9).Value = "sequence"
$newSheet.Cells.Item(2, 10).Value = "Net charge at pH 7"

$obRows = @(
    @("TRINITY_DN138200_c2_g1_i1.p2", [double]"-0.41388888888888797", [double]"6.4714511567774098"),
    @("TRINITY_DN141075_c3_g1_i1.p9", [double]"-1.69999999999999", [double]"3.4579619161958002"),
    @("TRINITY_DN143020_c3_g1_i1.p2", [double]"0.149999999999999", [double]"3.7611275943461902"),
    @("TRINITY_DN143683_c0_g2_i1.p14", [double]"-1.6666666666666701E-2", [double]"0.49935568654234302"),
    @("TRINITY_DN144897_c3_g1_i5.p10", [double]"-1.7222222222222201", [double]"1.5011571476165499"),
    @("TRINITY_DN385865_c1_g1_i1.p3", [double]"-0.79199999999999904", [double]"0.76313559624419802")
)

$r = 3
foreach ($row in $obRows) {
    $newSheet.Cells.Item($r, 7).Value = ">" + $row[0]
    $newSheet.Cells.Item($r, 8).Value = $row[1]
    $newSheet.Cells.Item($r, 9).Value = $row[0]
    $newSheet.Cells.Item($r, 10).Value = $row[2]
    $r++
}

# --- "Ob" rows of the combined table (columns A:D) ---
$r = 2
foreach ($row in $obRows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = "Ob"
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $r++
}

# --- "OMS" group source block (columns G:J) ---
$newSheet.Cells.Item(10, 7).Value = "sequence"
$newSheet.Cells.Item(10, 8).Value = "GRAVY index"
$newSheet.Cells.Item(10, 9).Value = "sequence"
$newSheet.Cells.Item(10, 10).Value = "Net charge at pH 7"

$omsRows = @(
    @("TRINITY_DN125192_c2_g1_i1.p8", [double]"0.27777777777777701", [double]"1.49934568664234"),
    @("TRINITY_DN133076_c0_g4_i1.p5", [double]"-0.82307692307692304", [double]"-0.32226295389458598"),
    @("TRINITY_DN138681_c7_g1_i1.p3", [double]"-0.75151515151515103", [double]"6.7489650825674898"),
    @("TRINITY_DN142115_c6_g1_i3.p7", [double]"-0.44444444444444398", [double]"-0.499523552517811"),
    @("TRINITY_DN142526_c5_g1_i7.p5", [double]"0.39999999999999902", [double]"1.7561056101932899"),
    @("TRINITY_DN142823_c1_g1_i2.p2", [double]"-0.32380952380952299", [double]"1.9260375560165699"),
    @("TRINITY_DN143096_c2_g1_i1.p6", [double]"-1.5538461538461501", [double]"5.7530986072962902"),
    @("TRINITY_DN144187_c0_g1_i10.p9", [double]"-3.3571428571428501", [double]"0.76003683310635395"),
    @("TRINITY_DN144483_c7_g1_i1.p5", [double]"0.628571428571428", [double]"-0.32729493794748299"),
    @("TRINITY_DN144624_c4_g1_i1.p5", [double]"0.55714285714285705", [double]"1.49934568664234"),
    @("TRINITY_DN144653_c7_g8_i1.p24", [double]"-2.2000000000000002", [double]"-4.4858548868306896"),
    @("TRINITY_DN145055_c3_g1_i3.p7", [double]"1.99999999999999E-2", [double]"-0.499523552517811"),
    @("TRINITY_DN145062_c3_g1_i3.p10", [double]"-3.73999999999999", [double]"1.83383061570099"),
    @("TRINITY_DN145245_c6_g1_i2.p2", [double]"-1.0344827586206799E-2", [double]"1.7584488922138199"),
    @("TRINITY_DN145391_c5_g1_i9.p5", [double]"-2.57777777777777", [double]"0.50396760958976194"),
    @("TRINITY_DN145575_c0_g1_i3.p4", [double]"-1.67777777777777", [double]"0.50215614861555302"),
    @("TRINITY_DN612198_c2_g1_i1.p1", [double]"-1.5414634146341399", [double]"3.9380451106218302"),
    @("TRINITY_DN7230_c1_g2_i2.p9", [double]"-1.0625", [double]"0.50297860849076204"),
    @("TRINITY_DN851511_c0_g1_i1.p12", [double]"1.3285714285714201", [double]"-0.49964531245865501")
)

$r = 11
foreach ($row in $omsRows) {
    $newSheet.Cells.Item($r, 7).Value = ">" + $row[0]
    $newSheet.Cells.Item($r, 8).Value = $row[1]
    $newSheet.Cells.Item($r, 9).Value = $row[0]
    $newSheet.Cells.Item($r, 10).Value = $row[2]
    $r++
}

# --- "OMS" rows of the combined table (columns A:D) ---
$r = 8
foreach ($row in $omsRows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = "OMS"
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $r++
}

[void]$newSheet.Range("A1:D26").Select()

# --- Selection change on the amps_overexpressed_info sheet ---
$infoSheet = $wb.Worksheets.Item("amps_overexpressed_info")
[void]$infoSheet.Activate()
[void]$infoSheet.Range("C16").Select()
